$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial value for every data row (rows 2-397).
# All of these were updated from 45189 (2023-09-20) to 45190 (2023-09-21).
$ws.Range("C2:C397").Value = 45190
